$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The course-name column (G) and meeting-days column (H) were swapped for
# each course's primary listing row. Fix it by swapping the values back.
$rows = @(2, 6, 10, 12, 14, 15, 16, 17, 18, 19, 20)

foreach ($r in $rows) {
    $gCell = $ws.Cells.Item($r, 7)   # column G
    $hCell = $ws.Cells.Item($r, 8)   # column H

    $gValue = $gCell.Value()
    $hValue = $hCell.Value()

    $gCell.Value = $hValue
    $hCell.Value = $gValue
}

# Leave the selection where the author's session ended up.
[void]$ws.Range("N2").Select()

